{"js": "// \"th\u00eam m\u00e1y t\u00ednh \u0111\u1ec3 b\u00e0n v\u00e0o csdl\"\n//\n// The document has a list entry \"M\u00e1y t\u00ednh b\u00e0n: 1\" (desktop computer count).\n// This edit:\n//   1. Strips the trailing count \"1\" from that line, leaving \"M\u00e1y t\u00ednh b\u00e0n: \".\n//   2. Adds a new sub-item right after it, one list level deeper, recording\n//      a specific brand/count: \"Lenovo: 1\".\n\n// Locate the unique paragraph holding the full phrase so we never touch any\n// of the many other \"...: 1\" / \"Lenovo: 1\" lines elsewhere in this (long)\n// device inventory document.\nconst matches = context.document.body.search(\"M\u00e1y t\u00ednh b\u00e0n: 1\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"M\u00e1y t\u00ednh b\u00e0n: 1\" in the document.');\n}\n\nconst hit = matches.items[0];\nconst para = hit.paragraphs.getFirst();\npara.load(\"text\");\nawait context.sync();\n\n// Replace only the trailing \": 1\" of this paragraph with \": \" (search is\n// scoped to the paragraph's own range, so it cannot match text elsewhere).\nconst tail = para.search(\": 1\", { matchCase: true });\ntail.load(\"items\");\nawait context.sync();\n\nconst tailHit = tail.items[tail.items.length - 1];\ntailHit.insertText(\": \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Insert the new \"Lenovo: 1\" line directly after, then push it one level\n// deeper in the same list (ilvl 4 -> 6).\nconst newPara = para.insertParagraph(\"Lenovo: 1\", Word.InsertLocation.after);\nawait context.sync();\n\nconst listItem = newPara.listItemOrNullObject;\nlistItem.load(\"level\");\nawait context.sync();\n\nif (!listItem.isNullObject) {\n  listItem.level = 6;\n  await context.sync();\n}\n", "ps1": "# \"th\u00eam m\u00e1y t\u00ednh \u0111\u1ec3 b\u00e0n v\u00e0o csdl\"\n#\n# The document has a list entry \"M\u00e1y t\u00ednh b\u00e0n: 1\" (desktop computer count).\n# This edit:\n#   1. Strips the trailing count \"1\" from that line, leaving \"M\u00e1y t\u00ednh b\u00e0n: \".\n#   2. Adds a new sub-item right after it, one list level deeper, recording\n#      a specific brand/count: \"Lenovo: 1\".\n\n$d = $word.ActiveDocument\n\n# Locate the unique paragraph by searching for the whole phrase, so we never\n# touch any of the many other \"...: 1\" / \"Lenovo: 1\" lines elsewhere in this\n# (long) device inventory document.\n$find = $d.Content\n$found = $find.Find.Execute(\"M\u00e1y t\u00ednh b\u00e0n: 1\")\nif (-not $found) {\n    throw 'Could not find \"M\u00e1y t\u00ednh b\u00e0n: 1\" in the document.'\n}\n\n$matchEnd = $find.End\n\n# Replace only the trailing \": 1\" (the last 3 characters of the match) with\n# \": \" -- delete it first, then insert the replacement via InsertAfter on a\n# collapsed range so the untouched \"b\u00e0n\" run ahead of it is left alone\n# instead of being merged with the edited text.\n$tail = $d.Range($matchEnd - 3, $matchEnd)\n$tail.Delete()\n$insertPoint = $d.Range($matchEnd - 3, $matchEnd - 3)\n$insertPoint.InsertAfter(\": \")\n\n# Add the new list line right after this paragraph, then push it one level\n# deeper in the same list (ilvl 4 -> 6, i.e. ListLevelNumber 5 -> 7).\n$para = $insertPoint.Paragraphs(1)\n$para.Range.InsertParagraphAfter()\n\n$newPara = $para.Next()\n$newPara.Range.InsertAfter(\"Lenovo: 1\")\n$newPara.Range.ListFormat.ListLevelNumber = 7\n"}
